$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New diary entry: "17 marras" (row 31) ---------------------------------
$ws.Range("A31").Value = "17 marras"
$ws.Range("C31").Value = "Oppikirjasta 213-228(kertaus), 231-"
$ws.Range("B31").Value = "10.30-11.00, 12.15-"
$ws.Range("D31").Value = "Rigidbody moottorikoodin kokonaisuuden kertausvilkuilu ( tässä kohtaa ei tehdä perässä uutta demoa näistä, vaan tehdään sitten kun on törmäykset), johdanto törmäyksiin kovilla kappaleilla, renderer/utility kirjaston päivitystä"

# Same look-and-feel as the rest of the log: wrapped text, "Kello" column
# formatted as a time, and a tall auto-fit-like row to match the long note.
$ws.Range("B31:D31").WrapText = $true
$ws.Range("B31").NumberFormat = "h:mm"
$ws.Rows.Item(31).RowHeight = 72.5

# The previous last row (30) picked up an extra formatted-but-empty cell
# in the "META" column, matching the other rows' wrap formatting.
$ws.Range("F30").WrapText = $true

# Leave the selection / scroll position where the user ended up after
# typing the new entry.
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select() | Out-Null
